# Update "想去人数" (want-to-go count) figures in column F for the two
# sheets that carry the full event listing: "展览" and "全部类型".
# The other two sheets ("演出" and "本地生活") only contain a header row
# and are left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F (applies to both sheets; both sheets
# converge on the same updated numbers even though sheet1/"展览" and
# sheet4/"全部类型" started from slightly different F3 values).
$updates = @{
    2  = 7747
    3  = 7577
    5  = 187
    8  = 124
    9  = 103
    10 = 145
    11 = 220
    12 = 683
    13 = 107
    14 = 1144
    16 = 41
    19 = 95
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
